$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sun Feb  5 13:50:20 UTC 2023 symbol-list refresh: updated Price (D) and
# Volume(1h) (E) columns for the affected rows. Cells are stored as literal
# text in the source workbook (inlineStr, no numeric style), so each cell is
# switched to the Text number format ("@") before its value is written -
# this keeps Excel from re-interpreting strings like "329.50" or "-0.41%"
# as a number/percentage. Formats are applied one cell at a time (not as a
# single unioned Range) since union-range NumberFormat assignment only
# reliably reaches the first area.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "329.50"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.41%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "43.21"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.71%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.592"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.75%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08202"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.61%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.765"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.51%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.423"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-2.31%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.917"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-5.89%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.847"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-4.84%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9429"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.38%"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-5.22%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1926"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.19%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09791"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "4.78%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.04491"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "14.80%"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.83%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001277"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.26%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006043"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.79%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.501"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.71%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.794"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "5.87%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.57%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2549"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "6.09%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04400"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.38%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001242"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.50%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004367"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.17%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001236"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "3.10%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004009"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "31.64%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02825"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-0.11%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05725"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.13%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007920"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.39%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009896"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "10.77%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.52%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-2.10%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009743"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-17.26%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007321"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "4.36%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.58%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003389"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "6.71%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.24%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.58%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.58%"
